$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rework the "Pid Mode" description (C4): reorder / reword the enumerated options
$newText = "0 -> Sensorwerte des Pyrometers und Leistungsmessers`n1 -> State des Buttons`n2 -> Status des Potentiometers`n3 -> Intern generierter Sägezahn"
$ws.Range("C4").Value = $newText

# Widen column C slightly to accommodate the revised text
$ws.Columns("C").ColumnWidth = 50.5

# Row 4 height should remain at 60pt (matches the now 4-line wrapped comment)
$ws.Rows("4").RowHeight = 60

# Leave the sheet's selection spanning the whole used range instead of a single stray cell
$ws.Range("A1:C4").Select() | Out-Null
